$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ row=25; A=24; E=45226.6875;          F="Kerala Blasters"; G=2; H="Odisha FC";  I=1; J=2.1;  K="21/10/2023 18:13"; L=2.3;  M="27/10/2023 16:28"; N=3.48; O="21/10/2023 18:13"; P=3.41; Q="27/10/2023 16:29"; R=3.43; S="21/10/2023 18:13"; T=3.16; U="27/10/2023 16:29"; V="https://www.betexplorer.com/football/india/isl/kerala-blasters-odisha-fc/UuzSNj1G/" },
    @{ row=26; A=25; E=45227.6875;          F="Mumbai City";     G=1; H="Hyderabad";   I=1; J=1.66; K="28/10/2023 13:33"; L=1.66; M="28/10/2023 13:33"; N=4.08; O="28/10/2023 13:33"; P=4.08; Q="28/10/2023 13:33"; R=4.94; S="28/10/2023 13:33"; T=4.94; U="28/10/2023 13:33"; V="https://www.betexplorer.com/football/india/isl/mumbai-city-hyderabad/rqvWMAGM/" },
    @{ row=27; A=26; E=45228.64583333334;   F="Chennaiyin";      G=5; H="Punjab";      I=1; J=1.72; K="29/10/2023 13:13"; L=1.82; M="29/10/2023 15:09"; N=3.86; O="29/10/2023 13:13"; P=3.85; Q="29/10/2023 15:09"; R=4.6;  S="29/10/2023 13:13"; T=4.19; U="29/10/2023 15:09"; V="https://www.betexplorer.com/football/india/isl/chennaiyin-fc-minerva-punjab/SGF2Hlvp/" }
)

foreach ($r in $rows) {
    $rowNum = $r.row

    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = "india"
    $ws.Cells.Item($rowNum, 3).Value = "isl"
    $ws.Cells.Item($rowNum, 4).Value = "2023-2024"
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
    $ws.Cells.Item($rowNum, 8).Value = $r.H
    $ws.Cells.Item($rowNum, 9).Value = $r.I
    $ws.Cells.Item($rowNum, 10).Value = $r.J
    $ws.Cells.Item($rowNum, 11).Value = $r.K
    $ws.Cells.Item($rowNum, 12).Value = $r.L
    $ws.Cells.Item($rowNum, 13).Value = $r.M
    $ws.Cells.Item($rowNum, 14).Value = $r.N
    $ws.Cells.Item($rowNum, 15).Value = $r.O
    $ws.Cells.Item($rowNum, 16).Value = $r.P
    $ws.Cells.Item($rowNum, 17).Value = $r.Q
    $ws.Cells.Item($rowNum, 18).Value = $r.R
    $ws.Cells.Item($rowNum, 19).Value = $r.S
    $ws.Cells.Item($rowNum, 20).Value = $r.T
    $ws.Cells.Item($rowNum, 21).Value = $r.U
    $ws.Cells.Item($rowNum, 22).Value = $r.V
}

# Copy cell formatting from the last existing data row (row 24) down to the
# newly added rows, matching the original style indices (s="1" on column A,
# s="2" on column E) without creating duplicate style entries.
$ws.Range("A24").Copy()
$ws.Range("A25:A27").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("E24").Copy()
$ws.Range("E25:E27").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

